# listeners, screenshot and more tests for login page added
#
# Adds two new worksheets ("errorInfoMessages" and "credentialsErrorMessages")
# after the existing "credentials" sheet, populates them with key/value
# lookup data (mirroring the style of the "credentials" sheet), and updates
# the selection/active-tab state so the new "credentialsErrorMessages" sheet
# ends up active.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. Add "errorInfoMessages" sheet (becomes sheetId 2 / rId2) right after the
#    "credentials" sheet.
# ---------------------------------------------------------------------------
$wsErrorInfo = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsErrorInfo.Name = "errorInfoMessages"

# Keys first (column A top-to-bottom), then values (column B top-to-bottom) -
# matches the shared-string insertion order of the authored workbook.
$wsErrorInfo.Range("A1").Value = "key"
$wsErrorInfo.Range("B1").Value = "errorMessage"
$wsErrorInfo.Range("A2").Value = "missingFirstName"
$wsErrorInfo.Range("A3").Value = "missingLastName"
$wsErrorInfo.Range("A4").Value = "missingZipCode"
$wsErrorInfo.Range("B2").Value = "Error: First Name is required"
$wsErrorInfo.Range("B3").Value = "Error: Last Name is required"
$wsErrorInfo.Range("B4").Value = "Error: Postal Code is required"

# Re-use the same header / body formatting as the "credentials" sheet.
$ws1.Range("A1:B1").Copy()
$wsErrorInfo.Range("A1:B1").PasteSpecial(-4122)
$ws1.Range("A2:B2").Copy()
$wsErrorInfo.Range("A2:B4").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$wsErrorInfo.Columns.Item(1).ColumnWidth = 25.7
$wsErrorInfo.Columns.Item(2).ColumnWidth = 40.7
[void]$wsErrorInfo.Range("A1:B4").Select()

# ---------------------------------------------------------------------------
# 2. Add "credentialsErrorMessages" sheet (becomes sheetId 3 / rId3) after
#    "errorInfoMessages".
# ---------------------------------------------------------------------------
$wsCredError = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsCredError.Name = "credentialsErrorMessages"

$wsCredError.Range("A1").Value = "key"
$wsCredError.Range("B1").Value = "message"
$wsCredError.Range("A2").Value = "invalid"
$wsCredError.Range("A3").Value = "locked"
$wsCredError.Range("B3").Value = "Sorry, this user has been locked out."
$wsCredError.Range("B2").Value = "Username and password do not match any user in this service."
$wsCredError.Range("B16").Value = " "

$ws1.Range("A1:B1").Copy()
$wsCredError.Range("A1:B1").PasteSpecial(-4122)
$ws1.Range("A2:B2").Copy()
$wsCredError.Range("A2:B3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$wsCredError.Columns.Item(1).ColumnWidth = 25.7
$wsCredError.Columns.Item(2).ColumnWidth = 89.9

# ---------------------------------------------------------------------------
# 3. Fix up row heights to mirror the "credentials" sheet's header / body
#    rows (26.25 for the header row, 23.25 for the data rows).
# ---------------------------------------------------------------------------
$wsErrorInfo.Rows.Item(1).RowHeight = 26.25
$wsErrorInfo.Rows.Item(2).RowHeight = 23.25
$wsErrorInfo.Rows.Item(3).RowHeight = 23.25
$wsErrorInfo.Rows.Item(4).RowHeight = 23.25

$wsCredError.Rows.Item(1).RowHeight = 26.25
$wsCredError.Rows.Item(2).RowHeight = 23.25
$wsCredError.Rows.Item(3).RowHeight = 23.25

# ---------------------------------------------------------------------------
# 4. Update selection on the original "credentials" sheet: it is no longer
#    the selected tab, and its selection becomes the full data range.
# ---------------------------------------------------------------------------
[void]$ws1.Activate()
[void]$ws1.Range("A1:B6").Select()

# ---------------------------------------------------------------------------
# 5. Final selection / active sheet: "credentialsErrorMessages" ends up
#    selected (tabSelected) with B14 as its active cell, and becomes the
#    workbook's active tab (activeTab="2").
# ---------------------------------------------------------------------------
[void]$wsCredError.Activate()
[void]$wsCredError.Range("B14").Select()
